# Mon_TI_2024_03_28 refresh: update the NBA player-prop table in Sheet1.
# Columns (header row 1): A=Equipe B=Joueur C=Statut D=Poste E=B2B F=5M G=15M
# H=Saison I=GP J=<20 K=20-29 L=30-39 M=40+ N=M-1 O=M-2 P=M-3 Q=M-4 R=M-5
# S=dom_ext T=delta_dom_ext U=delta_B2B V=nombre_de_B2B W=M+1_H_A X=M+1_team
# Y=M-1_vs Z=M-1_score AA=M-2_vs AB=M-2_score AC=M-3_vs AD=M-3_score
# AE=M+2_H_A AF=M+2_team AG=M+3_H_A AH=M+3_team AI=M+4_H_A AJ=M+4_team
# AK=M+5_H_A AL=M+5_team

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: now BOS / Jaylen Brown (previously ATL / Dejounte Murray)
$ws.Range("A2").Value = "BOS"
$ws.Range("B2").Value = "Jaylen Brown"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = 39.8
$ws.Range("G2").Value = 37.5
$ws.Range("I2").Value = 11
$ws.Range("J2").Value = 0
$ws.Range("M2").Value = 6
$ws.Range("N2").Value = 25
$ws.Range("O2").Value = "-"
$ws.Range("P2").Value = 43
$ws.Range("Q2").Value = 32
$ws.Range("R2").Value = 45
$ws.Range("S2").Value = "@"
$ws.Range("T2").Value = -2.2
$ws.Range("U2").Value = ""
$ws.Range("V2").Value = ""
$ws.Range("W2").Value = "@"
$ws.Range("X2").Value = "ATL"
$ws.Range("Z2").Value = 25
$ws.Range("AA2").Value = "vs"
$ws.Range("AB2").Value = 28
$ws.Range("AC2").Value = "vs"
$ws.Range("AD2").Value = 26
$ws.Range("AE2").Value = "@"
$ws.Range("AF2").Value = "NOP"
$ws.Range("AH2").Value = "CHA"
$ws.Range("AJ2").Value = "OKC"
$ws.Range("AK2").Value = "vs"
$ws.Range("AL2").Value = "SAC"

# Row 3: now ATL / Dejounte Murray (previously BOS / Jaylen Brown)
$ws.Range("A3").Value = "ATL"
$ws.Range("B3").Value = "Dejounte Murray"
$ws.Range("E3").Value = "O"
$ws.Range("F3").Value = 38.8
$ws.Range("G3").Value = 37.1
$ws.Range("I3").Value = 14
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 2
$ws.Range("N3").Value = 49
$ws.Range("O3").Value = 28
$ws.Range("P3").Value = 59
$ws.Range("Q3").Value = 42
$ws.Range("R3").Value = 16
$ws.Range("S3").Value = "vs"
$ws.Range("T3").Value = 0.5
$ws.Range("U3").Value = -8
$ws.Range("V3").Value = 11
$ws.Range("W3").Value = "vs"
$ws.Range("X3").Value = "BOS"
$ws.Range("AA3").Value = "@"
$ws.Range("AB3").Value = -3
$ws.Range("AE3").Value = "vs"
$ws.Range("AF3").Value = "MIL"
$ws.Range("AH3").Value = "CHI"
$ws.Range("AJ3").Value = "DET"
$ws.Range("AK3").Value = "@"
$ws.Range("AL3").Value = "DAL"

# Row 4: BOS / Kristaps Porzingis -- refreshed stat columns
$ws.Range("F4").Value = 26
$ws.Range("G4").Value = 35.1
$ws.Range("H4").Value = 32.8
$ws.Range("K4").Value = 2
$ws.Range("M4").Value = 1
$ws.Range("N4").Value = 24
$ws.Range("O4").Value = "-"
$ws.Range("P4").Value = 28
$ws.Range("Q4").Value = 15
$ws.Range("R4").Value = 33
$ws.Range("T4").Value = 0
$ws.Range("Y4").Value = "@"
$ws.Range("Z4").Value = 24
$ws.Range("AA4").Value = "vs"
$ws.Range("AB4").Value = 44

# Row 5: MIL / Khris Middleton -- refreshed stat columns + status
$ws.Range("C5").Value = "Probable"
$ws.Range("F5").Value = 25.8
$ws.Range("G5").Value = 25.4
$ws.Range("H5").Value = 23.8
$ws.Range("I5").Value = 4
$ws.Range("J5").Value = 1
$ws.Range("N5").Value = 17
$ws.Range("O5").Value = 23
$ws.Range("P5").Value = "-"
$ws.Range("Q5").Value = 41
$ws.Range("R5").Value = 36
$ws.Range("T5").Value = -0.6
